# Update Chirimoya price data rows (34-45 modified, 46-54 new) per weekly refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=34; D=44435; K="Cultivar IV Región"; L="Primera"; M=150; N=3200; O=3400; P=3300; Q="$/kilo (en caja de 15 kilos)"; R="Región de Coquimbo"; S=3300; T=1 },
    @{ Row=35; D=44435; K="Cultivar IV Región"; L="Segunda"; M=150; N=2800; O=3000; P=2900; Q="$/kilo (en caja de 15 kilos)"; R="Región de Coquimbo"; S=2900; T=1 },
    @{ Row=36; D=44438; K="Cultivar IV Región"; L="Especial"; M=80; N=24000; O=24000; P=24000; Q="$/bandeja 8 kilos"; R="Provincia del Elquí"; S=3000; T=8 },
    @{ Row=37; D=44438; K="Cultivar IV Región"; L="Extra (doble especial)"; M=100; N=3500; O=3500; P=3500; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=3500; T=1 },
    @{ Row=38; D=44438; K="Cultivar IV Región"; L="Extra (doble especial)"; M=85; N=3300; O=3300; P=3300; Q="$/kilo (en caja de 15 kilos)"; R="Región de Coquimbo"; S=3300; T=1 },
    @{ Row=39; D=44438; K="Cultivar IV Región"; L="Primera"; M=150; N=2900; O=2900; P=2900; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2900; T=1 },
    @{ Row=40; D=44438; K="Cultivar IV Región"; L="Primera"; M=95; N=3000; O=3000; P=3000; Q="$/kilo (en caja de 15 kilos)"; R="Región de Coquimbo"; S=3000; T=1 },
    @{ Row=41; D=44438; K="Cultivar IV Región"; L="Segunda"; M=100; N=2500; O=2500; P=2500; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=2500; T=1 },
    @{ Row=42; D=44438; K="Cultivar IV Región"; L="Segunda"; M=85; N=2600; O=2600; P=2600; Q="$/kilo (en caja de 15 kilos)"; R="Región de Coquimbo"; S=2600; T=1 },
    @{ Row=43; D=44159; K="Cultivar V Región"; L="Especial"; M=85; N=2000; O=2000; P=2000; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Quillota"; S=2000; T=1 },
    @{ Row=44; D=44159; K="Cultivar V Región"; L="Primera"; M=125; N=1700; O=1700; P=1700; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Quillota"; S=1700; T=1 },
    @{ Row=45; D=44159; K="Cultivar V Región"; L="Segunda"; M=95; N=1400; O=1400; P=1400; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Quillota"; S=1400; T=1 },
    @{ Row=46; D=44160; K="Cultivar IV Región"; L="Especial"; M=270; N=1900; O=1900; P=1900; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1900; T=1 },
    @{ Row=47; D=44160; K="Cultivar IV Región"; L="Primera"; M=270; N=1700; O=1700; P=1700; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1700; T=1 },
    @{ Row=48; D=44160; K="Cultivar IV Región"; L="Segunda"; M=270; N=1400; O=1400; P=1400; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1400; T=1 },
    @{ Row=49; D=44168; K="Cultivar IV Región"; L="Especial"; M=150; N=17000; O=17000; P=17000; Q="$/bandeja 8 kilos"; R="Provincia de Limarí"; S=2125; T=8 },
    @{ Row=50; D=44168; K="Cultivar IV Región"; L="Especial"; M=250; N=1800; O=1800; P=1800; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1800; T=1 },
    @{ Row=51; D=44168; K="Cultivar IV Región"; L="Primera"; M=350; N=15000; O=15000; P=15000; Q="$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1875; T=8 },
    @{ Row=52; D=44168; K="Cultivar IV Región"; L="Primera"; M=450; N=1500; O=1500; P=1500; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1500; T=1 },
    @{ Row=53; D=44168; K="Cultivar IV Región"; L="Segunda"; M=250; N=12000; O=12000; P=12000; Q="$/bandeja 8 kilos"; R="Provincia de Limarí"; S=1500; T=8 },
    @{ Row=54; D=44168; K="Cultivar IV Región"; L="Segunda"; M=350; N=1200; O=1200; P=1200; Q="$/kilo (en caja de 15 kilos)"; R="Provincia de Limarí"; S=1200; T=1 }
)

$constCols = @{
    A = 6
    B = "Mercado Mayorista Lo Valledor de Santiago"
    C = "Metropolitana"
    E = 13
    F = "Fruta"
    G = 100107
    H = "Otros"
    I = 100107002
    J = "Chirimoya"
}

foreach ($item in $rows) {
    $r = $item.Row

    foreach ($col in $constCols.Keys) {
        $ws.Cells.Item($r, [int][char]$col - [int][char]'A' + 1).Value = $constCols[$col]
    }

    $dCell = $ws.Range("D$r")
    $dCell.Value = $item.D
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("K$r").Value = $item.K
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("Q$r").Value = $item.Q
    $ws.Range("R$r").Value = $item.R
    $ws.Range("S$r").Value = $item.S
    $ws.Range("T$r").Value = $item.T
}
